$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Record 4 more problems finished in chapter 8: +1 in "String" row, +5 in "Linked L" row
$ws.Range("E8").Value = 7
$ws.Range("E9").Value = 5

# Move the active cell selection down to E10
$ws.Range("E10").Select()

$wb.Save()
